# Actualización de cambios: integracion Mercado Pago y Tickets CRM de Zoho.
# Replaces the demo "vino" catalog rows with the updated Mercado Pago style
# product list, clears the old numeric/stock data and repositions the
# stored selection/column widths to match the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = " text"
$ws.Range("B1").Value = "subText"
$ws.Range("C1").Value = "price"
$ws.Range("D1").Value = "stockDisponible"
$ws.Range("E1").Value = " image"

# Row 2
$ws.Range("A2").Value = "Single Vineyard"
$ws.Range("B2").Value = "Cabernet (x6)"
$ws.Range("C2").Value = 50
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

# Row 3
$ws.Range("A3").Value = "Almarada"
$ws.Range("B3").Value = "Malbec (x6)"
$ws.Range("C3").Value = 27000
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0

# Row 4
$ws.Range("A4").Value = "Núcleo3"
$ws.Range("B4").Value = "Malbec (x6)"
$ws.Range("C4").Value = 50
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0

# Row 5
$ws.Range("A5").Value = "Single Vineyard"
$ws.Range("B5").Value = "Chardonay (x6)"
$ws.Range("C5").Value = 50
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0

# Row 6
$ws.Range("A6").Value = "Núcleo3"
$ws.Range("B6").Value = "Malbec (x6)"
$ws.Range("C6").Value = 1000000
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0

# Row 7
$ws.Range("A7").Value = "Single Vineyard"
$ws.Range("B7").Value = "Malbec (x6)"
$ws.Range("C7").Value = 7000000
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0

# Column widths (closest achievable via character-based ColumnWidth,
# which snaps to the workbook's pixel grid)
$ws.Columns.Item(1).ColumnWidth = 15.25
$ws.Columns.Item(2).ColumnWidth = 14.25
$ws.Columns.Item(4).ColumnWidth = 22.6
$ws.Columns.Item(5).ColumnWidth = 37.7

# Restore the author's last-known selection
$ws.Range("E11").Select()
